$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (column H) so the shared style index is reused.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the new I/J columns with data for rows 2-13.
$values = @(
    @(1, 6),
    @(1, 6),
    @(2, 7),
    @(2, 6),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(8, 8),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}

$excel.CutCopyMode = $false
